$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Update column G width so stored width becomes 18 (from 17)
$ws.Columns.Item(7).ColumnWidth = 17.14

# Row 2
$ws.Cells.Item(2,1).Value = "2026-01-17"
$ws.Cells.Item(2,2).Value = "Booz"
$ws.Cells.Item(2,3).Value = "Yes"
$ws.Cells.Item(2,4).Value = "Nate Ament"
$ws.Cells.Item(2,5).Value = "TENN"
$ws.Cells.Item(2,6).Value = "UK@TENN"
$ws.Cells.Item(2,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(2,8).Value = 21
$ws.Cells.Item(2,9).Value = 17
$ws.Cells.Item(2,10).Value = 5
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 2
$ws.Cells.Item(2,13).Value = 0
$ws.Cells.Item(2,14).Value = 1
$ws.Cells.Item(2,15).Value = 35

# Row 3
$ws.Cells.Item(3,1).Value = "2026-01-17"
$ws.Cells.Item(3,2).Value = "Booz"
$ws.Cells.Item(3,3).Value = "Yes"
$ws.Cells.Item(3,4).Value = "Labaron Philon Jr."
$ws.Cells.Item(3,5).Value = "ALA"
$ws.Cells.Item(3,6).Value = "ALA@OU"
$ws.Cells.Item(3,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(3,8).Value = 3
$ws.Cells.Item(3,9).Value = 8
$ws.Cells.Item(3,10).Value = 3
$ws.Cells.Item(3,11).Value = 4
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0
$ws.Cells.Item(3,14).Value = 2
$ws.Cells.Item(3,15).Value = 21

# Row 4
$ws.Cells.Item(4,1).Value = "2026-01-17"
$ws.Cells.Item(4,2).Value = "Booz"
$ws.Cells.Item(4,3).Value = "Yes"
$ws.Cells.Item(4,4).Value = "Devin McGlockton"
$ws.Cells.Item(4,5).Value = "VAN"
$ws.Cells.Item(4,6).Value = "FLA@VAN"
$ws.Cells.Item(4,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(4,8).Value = -1
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(4,10).Value = 0
$ws.Cells.Item(4,11).Value = 0
$ws.Cells.Item(4,12).Value = 0
$ws.Cells.Item(4,13).Value = 0
$ws.Cells.Item(4,14).Value = 0
$ws.Cells.Item(4,15).Value = 6

# Row 5
$ws.Cells.Item(5,1).Value = "2026-01-17"
$ws.Cells.Item(5,2).Value = "Booz"
$ws.Cells.Item(5,3).Value = "No"
$ws.Cells.Item(5,4).Value = "Derrion Reid"
$ws.Cells.Item(5,5).Value = "OU"
$ws.Cells.Item(5,6).Value = "ALA@OU"
$ws.Cells.Item(5,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(5,8).Value = 19
$ws.Cells.Item(5,9).Value = 14
$ws.Cells.Item(5,10).Value = 8
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0
$ws.Cells.Item(5,13).Value = 0
$ws.Cells.Item(5,14).Value = 0
$ws.Cells.Item(5,15).Value = 24

# Row 6
$ws.Cells.Item(6,1).Value = "2026-01-17"
$ws.Cells.Item(6,2).Value = "CDL"
$ws.Cells.Item(6,3).Value = "Yes"
$ws.Cells.Item(6,4).Value = "Ja'Kobi Gillespie"
$ws.Cells.Item(6,5).Value = "TENN"
$ws.Cells.Item(6,6).Value = "UK@TENN"
$ws.Cells.Item(6,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(6,8).Value = 23
$ws.Cells.Item(6,9).Value = 24
$ws.Cells.Item(6,10).Value = 2
$ws.Cells.Item(6,11).Value = 8
$ws.Cells.Item(6,12).Value = 2
$ws.Cells.Item(6,13).Value = 0
$ws.Cells.Item(6,14).Value = 4
$ws.Cells.Item(6,15).Value = 33

# Row 7
$ws.Cells.Item(7,1).Value = "2026-01-17"
$ws.Cells.Item(7,2).Value = "CDL"
$ws.Cells.Item(7,3).Value = "Yes"
$ws.Cells.Item(7,4).Value = "Thomas Haugh"
$ws.Cells.Item(7,5).Value = "FLA"
$ws.Cells.Item(7,6).Value = "FLA@VAN"
$ws.Cells.Item(7,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(7,8).Value = 5
$ws.Cells.Item(7,9).Value = 4
$ws.Cells.Item(7,10).Value = 2
$ws.Cells.Item(7,11).Value = 0
$ws.Cells.Item(7,12).Value = 0
$ws.Cells.Item(7,13).Value = 0
$ws.Cells.Item(7,14).Value = 0
$ws.Cells.Item(7,15).Value = 7

# Row 8
$ws.Cells.Item(8,1).Value = "2026-01-17"
$ws.Cells.Item(8,2).Value = "CDL"
$ws.Cells.Item(8,3).Value = "No"
$ws.Cells.Item(8,4).Value = "Denzel Aberdeen"
$ws.Cells.Item(8,5).Value = "UK"
$ws.Cells.Item(8,6).Value = "UK@TENN"
$ws.Cells.Item(8,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(8,8).Value = 21
$ws.Cells.Item(8,9).Value = 22
$ws.Cells.Item(8,10).Value = 0
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 2
$ws.Cells.Item(8,13).Value = 0
$ws.Cells.Item(8,14).Value = 0
$ws.Cells.Item(8,15).Value = 29

# Row 9
$ws.Cells.Item(9,1).Value = "2026-01-17"
$ws.Cells.Item(9,2).Value = "CDL"
$ws.Cells.Item(9,3).Value = "No"
$ws.Cells.Item(9,4).Value = "Felix Okpara"
$ws.Cells.Item(9,5).Value = "TENN"
$ws.Cells.Item(9,6).Value = "UK@TENN"
$ws.Cells.Item(9,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(9,8).Value = 10
$ws.Cells.Item(9,9).Value = 5
$ws.Cells.Item(9,10).Value = 6
$ws.Cells.Item(9,11).Value = 0
$ws.Cells.Item(9,12).Value = 0
$ws.Cells.Item(9,13).Value = 1
$ws.Cells.Item(9,14).Value = 0
$ws.Cells.Item(9,15).Value = 26

# Row 10
$ws.Cells.Item(10,1).Value = "2026-01-17"
$ws.Cells.Item(10,2).Value = "CDL"
$ws.Cells.Item(10,3).Value = "No"
$ws.Cells.Item(10,4).Value = "Micah Handlogten"
$ws.Cells.Item(10,5).Value = "FLA"
$ws.Cells.Item(10,6).Value = "FLA@VAN"
$ws.Cells.Item(10,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(10,8).Value = 0
$ws.Cells.Item(10,9).Value = 0
$ws.Cells.Item(10,10).Value = 0
$ws.Cells.Item(10,11).Value = 0
$ws.Cells.Item(10,12).Value = 0
$ws.Cells.Item(10,13).Value = 0
$ws.Cells.Item(10,14).Value = 0
$ws.Cells.Item(10,15).Value = 2

# Row 11
$ws.Cells.Item(11,1).Value = "2026-01-17"
$ws.Cells.Item(11,2).Value = "CDL"
$ws.Cells.Item(11,3).Value = "No"
$ws.Cells.Item(11,4).Value = "Nijel Pack"
$ws.Cells.Item(11,5).Value = "OU"
$ws.Cells.Item(11,6).Value = "ALA@OU"
$ws.Cells.Item(11,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(11,8).Value = 0
$ws.Cells.Item(11,9).Value = 7
$ws.Cells.Item(11,10).Value = 1
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0
$ws.Cells.Item(11,13).Value = 0
$ws.Cells.Item(11,14).Value = 1
$ws.Cells.Item(11,15).Value = 17

# Row 12
$ws.Cells.Item(12,1).Value = "2026-01-17"
$ws.Cells.Item(12,2).Value = "Clay"
$ws.Cells.Item(12,3).Value = "Yes"
$ws.Cells.Item(12,4).Value = "Tyler Nickel"
$ws.Cells.Item(12,5).Value = "VAN"
$ws.Cells.Item(12,6).Value = "FLA@VAN"
$ws.Cells.Item(12,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(12,8).Value = 7
$ws.Cells.Item(12,9).Value = 5
$ws.Cells.Item(12,10).Value = 0
$ws.Cells.Item(12,11).Value = 0
$ws.Cells.Item(12,12).Value = 0
$ws.Cells.Item(12,13).Value = 2
$ws.Cells.Item(12,14).Value = 0
$ws.Cells.Item(12,15).Value = 6

# Row 13
$ws.Cells.Item(13,1).Value = "2026-01-17"
$ws.Cells.Item(13,2).Value = "Clay"
$ws.Cells.Item(13,3).Value = "Yes"
$ws.Cells.Item(13,4).Value = "Aden Holloway"
$ws.Cells.Item(13,5).Value = "ALA"
$ws.Cells.Item(13,6).Value = "ALA@OU"
$ws.Cells.Item(13,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(13,8).Value = 3
$ws.Cells.Item(13,9).Value = 7
$ws.Cells.Item(13,10).Value = 0
$ws.Cells.Item(13,11).Value = 1
$ws.Cells.Item(13,12).Value = 0
$ws.Cells.Item(13,13).Value = 0
$ws.Cells.Item(13,14).Value = 2
$ws.Cells.Item(13,15).Value = 18

# Row 14
$ws.Cells.Item(14,1).Value = "2026-01-17"
$ws.Cells.Item(14,2).Value = "Clay"
$ws.Cells.Item(14,3).Value = "No"
$ws.Cells.Item(14,4).Value = "Collin Chandler"
$ws.Cells.Item(14,5).Value = "UK"
$ws.Cells.Item(14,6).Value = "UK@TENN"
$ws.Cells.Item(14,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(14,8).Value = 17
$ws.Cells.Item(14,9).Value = 12
$ws.Cells.Item(14,10).Value = 3
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 3
$ws.Cells.Item(14,13).Value = 1
$ws.Cells.Item(14,14).Value = 1
$ws.Cells.Item(14,15).Value = 24

# Row 15
$ws.Cells.Item(15,1).Value = "2026-01-17"
$ws.Cells.Item(15,2).Value = "Clay"
$ws.Cells.Item(15,3).Value = "No"
$ws.Cells.Item(15,4).Value = "Jalil Bethea"
$ws.Cells.Item(15,5).Value = "ALA"
$ws.Cells.Item(15,6).Value = "ALA@OU"
$ws.Cells.Item(15,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(15,8).Value = 5
$ws.Cells.Item(15,9).Value = 4
$ws.Cells.Item(15,10).Value = 2
$ws.Cells.Item(15,11).Value = 0
$ws.Cells.Item(15,12).Value = 0
$ws.Cells.Item(15,13).Value = 0
$ws.Cells.Item(15,14).Value = 0
$ws.Cells.Item(15,15).Value = 6

# Row 16
$ws.Cells.Item(16,1).Value = "2026-01-17"
$ws.Cells.Item(16,2).Value = "Clay"
$ws.Cells.Item(16,3).Value = "No"
$ws.Cells.Item(16,4).Value = "Kirill Elatontsev"
$ws.Cells.Item(16,5).Value = "OU"
$ws.Cells.Item(16,6).Value = "ALA@OU"
$ws.Cells.Item(16,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(16,8).Value = 2
$ws.Cells.Item(16,9).Value = 0
$ws.Cells.Item(16,10).Value = 2
$ws.Cells.Item(16,11).Value = 1
$ws.Cells.Item(16,12).Value = 0
$ws.Cells.Item(16,13).Value = 0
$ws.Cells.Item(16,14).Value = 1
$ws.Cells.Item(16,15).Value = 10

# Row 17
$ws.Cells.Item(17,1).Value = "2026-01-17"
$ws.Cells.Item(17,2).Value = "Clay"
$ws.Cells.Item(17,3).Value = "No"
$ws.Cells.Item(17,4).Value = "Taylor Bol Bowen"
$ws.Cells.Item(17,5).Value = "ALA"
$ws.Cells.Item(17,6).Value = "ALA@OU"
$ws.Cells.Item(17,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(17,8).Value = 2
$ws.Cells.Item(17,9).Value = 0
$ws.Cells.Item(17,10).Value = 3
$ws.Cells.Item(17,11).Value = 0
$ws.Cells.Item(17,12).Value = 0
$ws.Cells.Item(17,13).Value = 1
$ws.Cells.Item(17,14).Value = 1
$ws.Cells.Item(17,15).Value = 12

# Row 18
$ws.Cells.Item(18,1).Value = "2026-01-17"
$ws.Cells.Item(18,2).Value = "Hal"
$ws.Cells.Item(18,3).Value = "Yes"
$ws.Cells.Item(18,4).Value = "Alex Condon"
$ws.Cells.Item(18,5).Value = "FLA"
$ws.Cells.Item(18,6).Value = "FLA@VAN"
$ws.Cells.Item(18,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(18,8).Value = 9
$ws.Cells.Item(18,9).Value = 9
$ws.Cells.Item(18,10).Value = 3
$ws.Cells.Item(18,11).Value = 0
$ws.Cells.Item(18,12).Value = 0
$ws.Cells.Item(18,13).Value = 1
$ws.Cells.Item(18,14).Value = 1
$ws.Cells.Item(18,15).Value = 7

# Row 19
$ws.Cells.Item(19,1).Value = "2026-01-17"
$ws.Cells.Item(19,2).Value = "Hal"
$ws.Cells.Item(19,3).Value = "Yes"
$ws.Cells.Item(19,4).Value = "Otega Oweh"
$ws.Cells.Item(19,5).Value = "UK"
$ws.Cells.Item(19,6).Value = "UK@TENN"
$ws.Cells.Item(19,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(19,8).Value = 3
$ws.Cells.Item(19,9).Value = 12
$ws.Cells.Item(19,10).Value = 2
$ws.Cells.Item(19,11).Value = 1
$ws.Cells.Item(19,12).Value = 2
$ws.Cells.Item(19,13).Value = 1
$ws.Cells.Item(19,14).Value = 3
$ws.Cells.Item(19,15).Value = 31

# Row 20
$ws.Cells.Item(20,1).Value = "2026-01-17"
$ws.Cells.Item(20,2).Value = "Hal"
$ws.Cells.Item(20,3).Value = "Yes"
$ws.Cells.Item(20,4).Value = "Boogie Fland"
$ws.Cells.Item(20,5).Value = "FLA"
$ws.Cells.Item(20,6).Value = "FLA@VAN"
$ws.Cells.Item(20,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(20,8).Value = 1
$ws.Cells.Item(20,9).Value = 2
$ws.Cells.Item(20,10).Value = 0
$ws.Cells.Item(20,11).Value = 0
$ws.Cells.Item(20,12).Value = 0
$ws.Cells.Item(20,13).Value = 0
$ws.Cells.Item(20,14).Value = 0
$ws.Cells.Item(20,15).Value = 7

# Row 21
$ws.Cells.Item(21,1).Value = "2026-01-17"
$ws.Cells.Item(21,2).Value = "Hal"
$ws.Cells.Item(21,3).Value = "No"
$ws.Cells.Item(21,4).Value = "Houston Mallette"
$ws.Cells.Item(21,5).Value = "ALA"
$ws.Cells.Item(21,6).Value = "ALA@OU"
$ws.Cells.Item(21,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(21,8).Value = 12
$ws.Cells.Item(21,9).Value = 8
$ws.Cells.Item(21,10).Value = 3
$ws.Cells.Item(21,11).Value = 1
$ws.Cells.Item(21,12).Value = 1
$ws.Cells.Item(21,13).Value = 0
$ws.Cells.Item(21,14).Value = 0
$ws.Cells.Item(21,15).Value = 17

# Row 22
$ws.Cells.Item(22,1).Value = "2026-01-17"
$ws.Cells.Item(22,2).Value = "Hal"
$ws.Cells.Item(22,3).Value = "No"
$ws.Cells.Item(22,4).Value = "Jadon Jones"
$ws.Cells.Item(22,5).Value = "OU"
$ws.Cells.Item(22,6).Value = "ALA@OU"
$ws.Cells.Item(22,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(22,8).Value = 9
$ws.Cells.Item(22,9).Value = 5
$ws.Cells.Item(22,10).Value = 2
$ws.Cells.Item(22,11).Value = 0
$ws.Cells.Item(22,12).Value = 3
$ws.Cells.Item(22,13).Value = 2
$ws.Cells.Item(22,14).Value = 0
$ws.Cells.Item(22,15).Value = 18

# Row 23
$ws.Cells.Item(23,1).Value = "2026-01-17"
$ws.Cells.Item(23,2).Value = "Hal"
$ws.Cells.Item(23,3).Value = "No"
$ws.Cells.Item(23,4).Value = "Jaylen Carey"
$ws.Cells.Item(23,5).Value = "TENN"
$ws.Cells.Item(23,6).Value = "UK@TENN"
$ws.Cells.Item(23,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(23,8).Value = 6
$ws.Cells.Item(23,9).Value = 6
$ws.Cells.Item(23,10).Value = 6
$ws.Cells.Item(23,11).Value = 3
$ws.Cells.Item(23,12).Value = 0
$ws.Cells.Item(23,13).Value = 0
$ws.Cells.Item(23,14).Value = 1
$ws.Cells.Item(23,15).Value = 21

# Row 24
$ws.Cells.Item(24,1).Value = "2026-01-17"
$ws.Cells.Item(24,2).Value = "Mark"
$ws.Cells.Item(24,3).Value = "Yes"
$ws.Cells.Item(24,4).Value = "Amari Allen"
$ws.Cells.Item(24,5).Value = "ALA"
$ws.Cells.Item(24,6).Value = "ALA@OU"
$ws.Cells.Item(24,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(24,8).Value = 12
$ws.Cells.Item(24,9).Value = 9
$ws.Cells.Item(24,10).Value = 5
$ws.Cells.Item(24,11).Value = 2
$ws.Cells.Item(24,12).Value = 1
$ws.Cells.Item(24,13).Value = 0
$ws.Cells.Item(24,14).Value = 1
$ws.Cells.Item(24,15).Value = 23

# Row 25
$ws.Cells.Item(25,1).Value = "2026-01-17"
$ws.Cells.Item(25,2).Value = "Mark"
$ws.Cells.Item(25,3).Value = "Yes"
$ws.Cells.Item(25,4).Value = "Xzayvier Brown"
$ws.Cells.Item(25,5).Value = "OU"
$ws.Cells.Item(25,6).Value = "ALA@OU"
$ws.Cells.Item(25,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(25,8).Value = 11
$ws.Cells.Item(25,9).Value = 11
$ws.Cells.Item(25,10).Value = 3
$ws.Cells.Item(25,11).Value = 3
$ws.Cells.Item(25,12).Value = 2
$ws.Cells.Item(25,13).Value = 0
$ws.Cells.Item(25,14).Value = 1
$ws.Cells.Item(25,15).Value = 24

# Row 26
$ws.Cells.Item(26,1).Value = "2026-01-17"
$ws.Cells.Item(26,2).Value = "Mark"
$ws.Cells.Item(26,3).Value = "Yes"
$ws.Cells.Item(26,4).Value = "Malachi Moreno"
$ws.Cells.Item(26,5).Value = "UK"
$ws.Cells.Item(26,6).Value = "UK@TENN"
$ws.Cells.Item(26,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(26,8).Value = 9
$ws.Cells.Item(26,9).Value = 3
$ws.Cells.Item(26,10).Value = 6
$ws.Cells.Item(26,11).Value = 2
$ws.Cells.Item(26,12).Value = 0
$ws.Cells.Item(26,13).Value = 3
$ws.Cells.Item(26,14).Value = 0
$ws.Cells.Item(26,15).Value = 22

# Row 27
$ws.Cells.Item(27,1).Value = "2026-01-17"
$ws.Cells.Item(27,2).Value = "Mark"
$ws.Cells.Item(27,3).Value = "Yes"
$ws.Cells.Item(27,4).Value = "Duke Miles"
$ws.Cells.Item(27,5).Value = "VAN"
$ws.Cells.Item(27,6).Value = "FLA@VAN"
$ws.Cells.Item(27,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(27,8).Value = 0
$ws.Cells.Item(27,9).Value = 0
$ws.Cells.Item(27,10).Value = 1
$ws.Cells.Item(27,11).Value = 2
$ws.Cells.Item(27,12).Value = 0
$ws.Cells.Item(27,13).Value = 0
$ws.Cells.Item(27,14).Value = 0
$ws.Cells.Item(27,15).Value = 6

# Row 28
$ws.Cells.Item(28,1).Value = "2026-01-17"
$ws.Cells.Item(28,2).Value = "Mark"
$ws.Cells.Item(28,3).Value = "No"
$ws.Cells.Item(28,4).Value = "J.P. Estrella"
$ws.Cells.Item(28,5).Value = "TENN"
$ws.Cells.Item(28,6).Value = "UK@TENN"
$ws.Cells.Item(28,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(28,8).Value = 5
$ws.Cells.Item(28,9).Value = 4
$ws.Cells.Item(28,10).Value = 0
$ws.Cells.Item(28,11).Value = 1
$ws.Cells.Item(28,12).Value = 1
$ws.Cells.Item(28,13).Value = 0
$ws.Cells.Item(28,14).Value = 1
$ws.Cells.Item(28,15).Value = 9

# Row 29
$ws.Cells.Item(29,1).Value = "2026-01-17"
$ws.Cells.Item(29,2).Value = "Mark"
$ws.Cells.Item(29,3).Value = "No"
$ws.Cells.Item(29,4).Value = "Xaivian Lee"
$ws.Cells.Item(29,5).Value = "FLA"
$ws.Cells.Item(29,6).Value = "FLA@VAN"
$ws.Cells.Item(29,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(29,8).Value = 3
$ws.Cells.Item(29,9).Value = 2
$ws.Cells.Item(29,10).Value = 1
$ws.Cells.Item(29,11).Value = 1
$ws.Cells.Item(29,12).Value = 0
$ws.Cells.Item(29,13).Value = 0
$ws.Cells.Item(29,14).Value = 0
$ws.Cells.Item(29,15).Value = 5

# Row 30
$ws.Cells.Item(30,1).Value = "2026-01-17"
$ws.Cells.Item(30,2).Value = "Ron"
$ws.Cells.Item(30,3).Value = "Yes"
$ws.Cells.Item(30,4).Value = "Mohamed Wague"
$ws.Cells.Item(30,5).Value = "OU"
$ws.Cells.Item(30,6).Value = "ALA@OU"
$ws.Cells.Item(30,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(30,8).Value = 8
$ws.Cells.Item(30,9).Value = 4
$ws.Cells.Item(30,10).Value = 5
$ws.Cells.Item(30,11).Value = 2
$ws.Cells.Item(30,12).Value = 1
$ws.Cells.Item(30,13).Value = 2
$ws.Cells.Item(30,14).Value = 1
$ws.Cells.Item(30,15).Value = 16

# Row 31
$ws.Cells.Item(31,1).Value = "2026-01-17"
$ws.Cells.Item(31,2).Value = "Ron"
$ws.Cells.Item(31,3).Value = "Yes"
$ws.Cells.Item(31,4).Value = "Rueben Chinyelu"
$ws.Cells.Item(31,5).Value = "FLA"
$ws.Cells.Item(31,6).Value = "FLA@VAN"
$ws.Cells.Item(31,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(31,8).Value = 3
$ws.Cells.Item(31,9).Value = 0
$ws.Cells.Item(31,10).Value = 2
$ws.Cells.Item(31,11).Value = 1
$ws.Cells.Item(31,12).Value = 0
$ws.Cells.Item(31,13).Value = 0
$ws.Cells.Item(31,14).Value = 0
$ws.Cells.Item(31,15).Value = 6

# Row 32
$ws.Cells.Item(32,1).Value = "2026-01-17"
$ws.Cells.Item(32,2).Value = "Ron"
$ws.Cells.Item(32,3).Value = "No"
$ws.Cells.Item(32,4).Value = "Mouhamed Dioubate"
$ws.Cells.Item(32,5).Value = "UK"
$ws.Cells.Item(32,6).Value = "UK@TENN"
$ws.Cells.Item(32,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(32,8).Value = 10
$ws.Cells.Item(32,9).Value = 10
$ws.Cells.Item(32,10).Value = 6
$ws.Cells.Item(32,11).Value = 0
$ws.Cells.Item(32,12).Value = 1
$ws.Cells.Item(32,13).Value = 0
$ws.Cells.Item(32,14).Value = 0
$ws.Cells.Item(32,15).Value = 27

# Row 33
$ws.Cells.Item(33,1).Value = "2026-01-17"
$ws.Cells.Item(33,2).Value = "Tar"
$ws.Cells.Item(33,3).Value = "Yes"
$ws.Cells.Item(33,4).Value = "Aiden Sherrell"
$ws.Cells.Item(33,5).Value = "ALA"
$ws.Cells.Item(33,6).Value = "ALA@OU"
$ws.Cells.Item(33,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(33,8).Value = 18
$ws.Cells.Item(33,9).Value = 14
$ws.Cells.Item(33,10).Value = 6
$ws.Cells.Item(33,11).Value = 1
$ws.Cells.Item(33,12).Value = 1
$ws.Cells.Item(33,13).Value = 1
$ws.Cells.Item(33,14).Value = 1
$ws.Cells.Item(33,15).Value = 21

# Row 34
$ws.Cells.Item(34,1).Value = "2026-01-17"
$ws.Cells.Item(34,2).Value = "Tar"
$ws.Cells.Item(34,3).Value = "Yes"
$ws.Cells.Item(34,4).Value = "Tae Davis"
$ws.Cells.Item(34,5).Value = "OU"
$ws.Cells.Item(34,6).Value = "ALA@OU"
$ws.Cells.Item(34,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(34,8).Value = 12
$ws.Cells.Item(34,9).Value = 14
$ws.Cells.Item(34,10).Value = 4
$ws.Cells.Item(34,11).Value = 0
$ws.Cells.Item(34,12).Value = 0
$ws.Cells.Item(34,13).Value = 0
$ws.Cells.Item(34,14).Value = 2
$ws.Cells.Item(34,15).Value = 24

# Row 35
$ws.Cells.Item(35,1).Value = "2026-01-17"
$ws.Cells.Item(35,2).Value = "Tar"
$ws.Cells.Item(35,3).Value = "Yes"
$ws.Cells.Item(35,4).Value = "Tyler Tanner"
$ws.Cells.Item(35,5).Value = "VAN"
$ws.Cells.Item(35,6).Value = "FLA@VAN"
$ws.Cells.Item(35,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(35,8).Value = 7
$ws.Cells.Item(35,9).Value = 3
$ws.Cells.Item(35,10).Value = 3
$ws.Cells.Item(35,11).Value = 0
$ws.Cells.Item(35,12).Value = 1
$ws.Cells.Item(35,13).Value = 1
$ws.Cells.Item(35,14).Value = 0
$ws.Cells.Item(35,15).Value = 8

# Row 36
$ws.Cells.Item(36,1).Value = "2026-01-17"
$ws.Cells.Item(36,2).Value = "Tar"
$ws.Cells.Item(36,3).Value = "No"
$ws.Cells.Item(36,4).Value = "Jalen Washington"
$ws.Cells.Item(36,5).Value = "VAN"
$ws.Cells.Item(36,6).Value = "FLA@VAN"
$ws.Cells.Item(36,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(36,8).Value = -2
$ws.Cells.Item(36,9).Value = 0
$ws.Cells.Item(36,10).Value = 0
$ws.Cells.Item(36,11).Value = 0
$ws.Cells.Item(36,12).Value = 0
$ws.Cells.Item(36,13).Value = 0
$ws.Cells.Item(36,14).Value = 0
$ws.Cells.Item(36,15).Value = 4

# Row 37
$ws.Cells.Item(37,1).Value = "2026-01-17"
$ws.Cells.Item(37,2).Value = "Undrafted"
$ws.Cells.Item(37,3).Value = "No"
$ws.Cells.Item(37,4).Value = "Jasper Johnson"
$ws.Cells.Item(37,5).Value = "UK"
$ws.Cells.Item(37,6).Value = "UK@TENN"
$ws.Cells.Item(37,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(37,8).Value = 16
$ws.Cells.Item(37,9).Value = 12
$ws.Cells.Item(37,10).Value = 2
$ws.Cells.Item(37,11).Value = 4
$ws.Cells.Item(37,12).Value = 0
$ws.Cells.Item(37,13).Value = 0
$ws.Cells.Item(37,14).Value = 1
$ws.Cells.Item(37,15).Value = 19

# Row 38
$ws.Cells.Item(38,1).Value = "2026-01-17"
$ws.Cells.Item(38,2).Value = "Undrafted"
$ws.Cells.Item(38,3).Value = "No"
$ws.Cells.Item(38,4).Value = "DeWayne Brown II"
$ws.Cells.Item(38,5).Value = "TENN"
$ws.Cells.Item(38,6).Value = "UK@TENN"
$ws.Cells.Item(38,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(38,8).Value = 11
$ws.Cells.Item(38,9).Value = 7
$ws.Cells.Item(38,10).Value = 2
$ws.Cells.Item(38,11).Value = 2
$ws.Cells.Item(38,12).Value = 1
$ws.Cells.Item(38,13).Value = 2
$ws.Cells.Item(38,14).Value = 1
$ws.Cells.Item(38,15).Value = 19

# Row 39
$ws.Cells.Item(39,1).Value = "2026-01-17"
$ws.Cells.Item(39,2).Value = "Undrafted"
$ws.Cells.Item(39,3).Value = "No"
$ws.Cells.Item(39,4).Value = "Amari Evans"
$ws.Cells.Item(39,5).Value = "TENN"
$ws.Cells.Item(39,6).Value = "UK@TENN"
$ws.Cells.Item(39,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(39,8).Value = 10
$ws.Cells.Item(39,9).Value = 8
$ws.Cells.Item(39,10).Value = 3
$ws.Cells.Item(39,11).Value = 0
$ws.Cells.Item(39,12).Value = 0
$ws.Cells.Item(39,13).Value = 1
$ws.Cells.Item(39,14).Value = 0
$ws.Cells.Item(39,15).Value = 16

# Row 40
$ws.Cells.Item(40,1).Value = "2026-01-17"
$ws.Cells.Item(40,2).Value = "Undrafted"
$ws.Cells.Item(40,3).Value = "No"
$ws.Cells.Item(40,4).Value = "AK Okereke"
$ws.Cells.Item(40,5).Value = "VAN"
$ws.Cells.Item(40,6).Value = "FLA@VAN"
$ws.Cells.Item(40,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(40,8).Value = 7
$ws.Cells.Item(40,9).Value = 7
$ws.Cells.Item(40,10).Value = 0
$ws.Cells.Item(40,11).Value = 0
$ws.Cells.Item(40,12).Value = 0
$ws.Cells.Item(40,13).Value = 0
$ws.Cells.Item(40,14).Value = 0
$ws.Cells.Item(40,15).Value = 6

# Row 41
$ws.Cells.Item(41,1).Value = "2026-01-17"
$ws.Cells.Item(41,2).Value = "Undrafted"
$ws.Cells.Item(41,3).Value = "No"
$ws.Cells.Item(41,4).Value = "Andrija Jelavić"
$ws.Cells.Item(41,5).Value = "UK"
$ws.Cells.Item(41,6).Value = "UK@TENN"
$ws.Cells.Item(41,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(41,8).Value = 5
$ws.Cells.Item(41,9).Value = 5
$ws.Cells.Item(41,10).Value = 1
$ws.Cells.Item(41,11).Value = 1
$ws.Cells.Item(41,12).Value = 0
$ws.Cells.Item(41,13).Value = 0
$ws.Cells.Item(41,14).Value = 0
$ws.Cells.Item(41,15).Value = 12

# Row 42
$ws.Cells.Item(42,1).Value = "2026-01-17"
$ws.Cells.Item(42,2).Value = "Undrafted"
$ws.Cells.Item(42,3).Value = "No"
$ws.Cells.Item(42,4).Value = "Brandon Garrison"
$ws.Cells.Item(42,5).Value = "UK"
$ws.Cells.Item(42,6).Value = "UK@TENN"
$ws.Cells.Item(42,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(42,8).Value = 4
$ws.Cells.Item(42,9).Value = 2
$ws.Cells.Item(42,10).Value = 2
$ws.Cells.Item(42,11).Value = 1
$ws.Cells.Item(42,12).Value = 0
$ws.Cells.Item(42,13).Value = 0
$ws.Cells.Item(42,14).Value = 0
$ws.Cells.Item(42,15).Value = 18

# Row 43
$ws.Cells.Item(43,1).Value = "2026-01-17"
$ws.Cells.Item(43,2).Value = "Undrafted"
$ws.Cells.Item(43,3).Value = "No"
$ws.Cells.Item(43,4).Value = "Isaiah Brown"
$ws.Cells.Item(43,5).Value = "FLA"
$ws.Cells.Item(43,6).Value = "FLA@VAN"
$ws.Cells.Item(43,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(43,8).Value = 3
$ws.Cells.Item(43,9).Value = 2
$ws.Cells.Item(43,10).Value = 1
$ws.Cells.Item(43,11).Value = 0
$ws.Cells.Item(43,12).Value = 0
$ws.Cells.Item(43,13).Value = 0
$ws.Cells.Item(43,14).Value = 0
$ws.Cells.Item(43,15).Value = 2

# Row 44
$ws.Cells.Item(44,1).Value = "2026-01-17"
$ws.Cells.Item(44,2).Value = "Undrafted"
$ws.Cells.Item(44,3).Value = "No"
$ws.Cells.Item(44,4).Value = "Amaree Abram"
$ws.Cells.Item(44,5).Value = "TENN"
$ws.Cells.Item(44,6).Value = "UK@TENN"
$ws.Cells.Item(44,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(44,8).Value = 2
$ws.Cells.Item(44,9).Value = 3
$ws.Cells.Item(44,10).Value = 1
$ws.Cells.Item(44,11).Value = 0
$ws.Cells.Item(44,12).Value = 0
$ws.Cells.Item(44,13).Value = 0
$ws.Cells.Item(44,14).Value = 1
$ws.Cells.Item(44,15).Value = 5

# Row 45
$ws.Cells.Item(45,1).Value = "2026-01-17"
$ws.Cells.Item(45,2).Value = "Undrafted"
$ws.Cells.Item(45,3).Value = "No"
$ws.Cells.Item(45,4).Value = "Bishop Boswell"
$ws.Cells.Item(45,5).Value = "TENN"
$ws.Cells.Item(45,6).Value = "UK@TENN"
$ws.Cells.Item(45,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(45,8).Value = 2
$ws.Cells.Item(45,9).Value = 4
$ws.Cells.Item(45,10).Value = 3
$ws.Cells.Item(45,11).Value = 1
$ws.Cells.Item(45,12).Value = 0
$ws.Cells.Item(45,13).Value = 0
$ws.Cells.Item(45,14).Value = 2
$ws.Cells.Item(45,15).Value = 26

# Row 46
$ws.Cells.Item(46,1).Value = "2026-01-17"
$ws.Cells.Item(46,2).Value = "Undrafted"
$ws.Cells.Item(46,3).Value = "No"
$ws.Cells.Item(46,4).Value = "Tyler Harris"
$ws.Cells.Item(46,5).Value = "VAN"
$ws.Cells.Item(46,6).Value = "FLA@VAN"
$ws.Cells.Item(46,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(46,8).Value = 1
$ws.Cells.Item(46,9).Value = 2
$ws.Cells.Item(46,10).Value = 0
$ws.Cells.Item(46,11).Value = 0
$ws.Cells.Item(46,12).Value = 0
$ws.Cells.Item(46,13).Value = 0
$ws.Cells.Item(46,14).Value = 0
$ws.Cells.Item(46,15).Value = 4

# Row 47
$ws.Cells.Item(47,1).Value = "2026-01-17"
$ws.Cells.Item(47,2).Value = "Undrafted"
$ws.Cells.Item(47,3).Value = "No"
$ws.Cells.Item(47,4).Value = "Urban Klavzar"
$ws.Cells.Item(47,5).Value = "FLA"
$ws.Cells.Item(47,6).Value = "FLA@VAN"
$ws.Cells.Item(47,7).Value = "12:04 - 1st Half"
$ws.Cells.Item(47,8).Value = 1
$ws.Cells.Item(47,9).Value = 0
$ws.Cells.Item(47,10).Value = 0
$ws.Cells.Item(47,11).Value = 1
$ws.Cells.Item(47,12).Value = 0
$ws.Cells.Item(47,13).Value = 0
$ws.Cells.Item(47,14).Value = 0
$ws.Cells.Item(47,15).Value = 4

# Row 48
$ws.Cells.Item(48,1).Value = "2026-01-17"
$ws.Cells.Item(48,2).Value = "Undrafted"
$ws.Cells.Item(48,3).Value = "No"
$ws.Cells.Item(48,4).Value = "Kam Williams"
$ws.Cells.Item(48,5).Value = "UK"
$ws.Cells.Item(48,6).Value = "UK@TENN"
$ws.Cells.Item(48,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(48,8).Value = 0
$ws.Cells.Item(48,9).Value = 2
$ws.Cells.Item(48,10).Value = 2
$ws.Cells.Item(48,11).Value = 1
$ws.Cells.Item(48,12).Value = 0
$ws.Cells.Item(48,13).Value = 0
$ws.Cells.Item(48,14).Value = 2
$ws.Cells.Item(48,15).Value = 15

# Row 49
$ws.Cells.Item(49,1).Value = "2026-01-17"
$ws.Cells.Item(49,2).Value = "Undrafted"
$ws.Cells.Item(49,3).Value = "No"
$ws.Cells.Item(49,4).Value = "London Jemison"
$ws.Cells.Item(49,5).Value = "ALA"
$ws.Cells.Item(49,6).Value = "ALA@OU"
$ws.Cells.Item(49,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(49,8).Value = 0
$ws.Cells.Item(49,9).Value = 2
$ws.Cells.Item(49,10).Value = 1
$ws.Cells.Item(49,11).Value = 0
$ws.Cells.Item(49,12).Value = 0
$ws.Cells.Item(49,13).Value = 1
$ws.Cells.Item(49,14).Value = 0
$ws.Cells.Item(49,15).Value = 10

# Row 50
$ws.Cells.Item(50,1).Value = "2026-01-17"
$ws.Cells.Item(50,2).Value = "Undrafted"
$ws.Cells.Item(50,3).Value = "No"
$ws.Cells.Item(50,4).Value = "Noah Williamson"
$ws.Cells.Item(50,5).Value = "ALA"
$ws.Cells.Item(50,6).Value = "ALA@OU"
$ws.Cells.Item(50,7).Value = "13:25 - 2nd Half"
$ws.Cells.Item(50,8).Value = 0
$ws.Cells.Item(50,9).Value = 0
$ws.Cells.Item(50,10).Value = 0
$ws.Cells.Item(50,11).Value = 0
$ws.Cells.Item(50,12).Value = 0
$ws.Cells.Item(50,13).Value = 0
$ws.Cells.Item(50,14).Value = 0
$ws.Cells.Item(50,15).Value = 3

# Row 51
$ws.Cells.Item(51,1).Value = "2026-01-17"
$ws.Cells.Item(51,2).Value = "Undrafted"
$ws.Cells.Item(51,3).Value = "No"
$ws.Cells.Item(51,4).Value = "Ethan Burg"
$ws.Cells.Item(51,5).Value = "TENN"
$ws.Cells.Item(51,6).Value = "UK@TENN"
$ws.Cells.Item(51,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(51,8).Value = -1
$ws.Cells.Item(51,9).Value = 0
$ws.Cells.Item(51,10).Value = 0
$ws.Cells.Item(51,11).Value = 0
$ws.Cells.Item(51,12).Value = 0
$ws.Cells.Item(51,13).Value = 0
$ws.Cells.Item(51,14).Value = 1
$ws.Cells.Item(51,15).Value = 3

# Row 52
$ws.Cells.Item(52,1).Value = "2026-01-17"
$ws.Cells.Item(52,2).Value = "Undrafted"
$ws.Cells.Item(52,3).Value = "No"
$ws.Cells.Item(52,4).Value = "Trent Noah"
$ws.Cells.Item(52,5).Value = "UK"
$ws.Cells.Item(52,6).Value = "UK@TENN"
$ws.Cells.Item(52,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(52,8).Value = -1
$ws.Cells.Item(52,9).Value = 0
$ws.Cells.Item(52,10).Value = 0
$ws.Cells.Item(52,11).Value = 0
$ws.Cells.Item(52,12).Value = 0
$ws.Cells.Item(52,13).Value = 0
$ws.Cells.Item(52,14).Value = 1
$ws.Cells.Item(52,15).Value = 3

# Row 53
$ws.Cells.Item(53,1).Value = "2026-01-17"
$ws.Cells.Item(53,2).Value = "Undrafted"
$ws.Cells.Item(53,3).Value = "No"
$ws.Cells.Item(53,4).Value = "Troy Henderson"
$ws.Cells.Item(53,5).Value = "TENN"
$ws.Cells.Item(53,6).Value = "UK@TENN"
$ws.Cells.Item(53,7).Value = "0:02 - 2nd Half"
$ws.Cells.Item(53,8).Value = -2
$ws.Cells.Item(53,9).Value = 0
$ws.Cells.Item(53,10).Value = 0
$ws.Cells.Item(53,11).Value = 1
$ws.Cells.Item(53,12).Value = 0
$ws.Cells.Item(53,13).Value = 0
$ws.Cells.Item(53,14).Value = 0
$ws.Cells.Item(53,15).Value = 6

$ws2 = $wb.Worksheets.Item("OwnerTotals")

# OwnerTotals row 2
$ws2.Cells.Item(2,1).Value = "Tar"
$ws2.Cells.Item(2,2).Value = 37
$ws2.Cells.Item(2,3).Value = 3

# OwnerTotals row 3
$ws2.Cells.Item(3,1).Value = "Mark"
$ws2.Cells.Item(3,2).Value = 32
$ws2.Cells.Item(3,3).Value = 4

# OwnerTotals row 4
$ws2.Cells.Item(4,1).Value = "CDL"
$ws2.Cells.Item(4,2).Value = 28
$ws2.Cells.Item(4,3).Value = 2

# OwnerTotals row 5
$ws2.Cells.Item(5,1).Value = "Booz"
$ws2.Cells.Item(5,2).Value = 23
$ws2.Cells.Item(5,3).Value = 3

# OwnerTotals row 6
$ws2.Cells.Item(6,1).Value = "Hal"
$ws2.Cells.Item(6,2).Value = 13
$ws2.Cells.Item(6,3).Value = 3

# OwnerTotals row 7
$ws2.Cells.Item(7,1).Value = "Ron"
$ws2.Cells.Item(7,2).Value = 11
$ws2.Cells.Item(7,3).Value = 2

# OwnerTotals row 8
$ws2.Cells.Item(8,1).Value = "Clay"
$ws2.Cells.Item(8,2).Value = 10
$ws2.Cells.Item(8,3).Value = 2

